$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1127.25
$ws.Range("I28").Value = 1127.25
$ws.Range("K28").Value = 1127.25
$ws.Range("M28").Value = -642.25
$ws.Range("H101").Value = 3021.889
$ws.Range("I101").Value = 601.6
$ws.Range("K101").Value = 1804.8
$ws.Range("M101").Value = -182.8000000000002
$ws.Range("H129").Value = 1846.5883
$ws.Range("I129").Value = 675.1667
$ws.Range("K129").Value = 2025.5001
$ws.Range("M129").Value = 2974.4999
$ws.Range("H132").Value = 12347.63
$ws.Range("I132").Value = 2498.0908
$ws.Range("K132").Value = 7494.2724
$ws.Range("M132").Value = -4964.2724
$ws.Range("H139").Value = 112629.836
$ws.Range("J139").Value = 112629.836
$ws.Range("L139").Value = 112629.836
$ws.Range("N139").Value = -122909.836
$ws.Range("H140").Value = 56265.875
$ws.Range("J140").Value = 54202.57
$ws.Range("L140").Value = 54202.57
$ws.Range("N140").Value = -64562.57
$ws.Range("H141").Value = 6411.731
$ws.Range("I141").Value = 5291.7827
$ws.Range("J141").Value = 14998
$ws.Range("K141").Value = 15875.3481
$ws.Range("L141").Value = 44994
$ws.Range("M141").Value = -10695.3481
$ws.Range("N141").Value = -55354

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 3991.6667
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H61").Value = 6765.433
$ws.Range("I61").Value = 7340.3335
$ws.Range("K61").Value = 7340.3335
$ws.Range("M61").Value = -7128.3335
$ws.Range("H74").Value = 1418.1538
$ws.Range("I74").Value = 1091.8667
$ws.Range("K74").Value = 1091.8667
$ws.Range("M74").Value = -217.8667
$ws.Range("H77").Value = 1418.1538
$ws.Range("I77").Value = 1091.8667
$ws.Range("K77").Value = 5459.333500000001
$ws.Range("M77").Value = -1091.333500000001
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H132").Value = 13227.982
$ws.Range("I132").Value = 14743.571
$ws.Range("J132").Value = 8984.333000000001
$ws.Range("K132").Value = 44230.713
$ws.Range("L132").Value = 26952.999
$ws.Range("M132").Value = -41700.713
$ws.Range("N132").Value = -32012.999
$ws.Range("H136").Value = 6765.433
$ws.Range("I136").Value = 7340.3335
$ws.Range("K136").Value = 22021.0005
$ws.Range("M136").Value = -19471.0005
$ws.Range("H140").Value = 114600
$ws.Range("J140").Value = 114600
$ws.Range("L140").Value = 114600
$ws.Range("N140").Value = -124960

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 585.1429000000001
$ws.Range("I22").Value = 585.1429000000001
$ws.Range("K22").Value = 585.1429000000001
$ws.Range("M22").Value = -412.1429000000001
$ws.Range("H134").Value = 2023.6666
$ws.Range("I134").Value = 1679.8036
$ws.Range("K134").Value = 5039.4108
$ws.Range("M134").Value = -2504.4108
$ws.Range("H140").Value = 99007.836
$ws.Range("J140").Value = 99007.836
$ws.Range("L140").Value = 99007.836
$ws.Range("N140").Value = -109367.836

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2344.5806
$ws.Range("I31").Value = 1133.6111
$ws.Range("J31").Value = 4021.3076
$ws.Range("K31").Value = 1133.6111
$ws.Range("L31").Value = 4021.3076
$ws.Range("M31").Value = -838.6111000000001
$ws.Range("N31").Value = -4611.3076
$ws.Range("H34").Value = 2344.5806
$ws.Range("I34").Value = 1133.6111
$ws.Range("J34").Value = 4021.3076
$ws.Range("K34").Value = 1133.6111
$ws.Range("L34").Value = 4021.3076
$ws.Range("M34").Value = -931.6111000000001
$ws.Range("N34").Value = -4425.3076
$ws.Range("H132").Value = 2149.5789
$ws.Range("I132").Value = 1935.7778
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 5807.3334
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -3277.3334
$ws.Range("N132").Value = -23054
$ws.Range("H141").Value = 113721.22
$ws.Range("J141").Value = 122999.125
$ws.Range("L141").Value = 122999.125
$ws.Range("N141").Value = -133359.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 246
$ws.Range("I8").Value = 246
$ws.Range("K8").Value = 738
$ws.Range("M8").Value = -599
$ws.Range("H12").Value = 365.53333
$ws.Range("J12").Value = 366.81818
$ws.Range("L12").Value = 1100.45454
$ws.Range("N12").Value = -1446.45454
$ws.Range("H92").Value = 706.0909
$ws.Range("J92").Value = 667.5
$ws.Range("L92").Value = 2002.5
$ws.Range("N92").Value = -4498.5
$ws.Range("H128").Value = 500000
$ws.Range("I128").Value = 500000
$ws.Range("K128").Value = 1500000
$ws.Range("M128").Value = -1495020
$ws.Range("H140").Value = 4339.8887
$ws.Range("I140").Value = 4339.8887
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 13019.6661
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -7839.666100000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3973862.8
$ws.Range("I70").Value = 4767635.5
$ws.Range("K70").Value = 4767635.5
$ws.Range("M70").Value = -4767365.5
$ws.Range("H73").Value = 3973862.8
$ws.Range("I73").Value = 4767635.5
$ws.Range("K73").Value = 4767635.5
$ws.Range("M73").Value = -4766699.5
$ws.Range("H102").Value = 9780.579
$ws.Range("I102").Value = 9559.357
$ws.Range("K102").Value = 9559.357
$ws.Range("M102").Value = -7937.357
$ws.Range("H107").Value = 9524570
$ws.Range("J107").Value = 933
$ws.Range("L107").Value = 933
$ws.Range("N107").Value = -4773
$ws.Range("H126").Value = 3699.7058
$ws.Range("J126").Value = 5849
$ws.Range("L126").Value = 17547
$ws.Range("N126").Value = -22487
$ws.Range("H132").Value = 3710.8948
$ws.Range("I132").Value = 3538.0625
$ws.Range("K132").Value = 10614.1875
$ws.Range("M132").Value = -8084.1875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5833.4287
$ws.Range("I40").Value = 5341.5454
$ws.Range("K40").Value = 5341.5454
$ws.Range("M40").Value = -5205.5454
$ws.Range("H46").Value = 6474.3105
$ws.Range("J46").Value = 6982.8696
$ws.Range("L46").Value = 6982.8696
$ws.Range("N46").Value = -7358.8696
$ws.Range("H132").Value = 3635.2856
$ws.Range("I132").Value = 2841.1924
$ws.Range("J132").Value = 5929.3335
$ws.Range("K132").Value = 8523.5772
$ws.Range("L132").Value = 17788.0005
$ws.Range("M132").Value = -5993.5772
$ws.Range("N132").Value = -22848.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10290252
$ws.Range("I132").Value = 1426109.5
$ws.Range("K132").Value = 4278328.5
$ws.Range("M132").Value = -4275798.5
$ws.Range("H136").Value = 8291.216
$ws.Range("I136").Value = 4738.643
$ws.Range("J136").Value = 8963.324000000001
$ws.Range("K136").Value = 14215.929
$ws.Range("L136").Value = 26889.972
$ws.Range("M136").Value = -11665.929
$ws.Range("N136").Value = -31989.972
$ws.Range("H138").Value = 71990
$ws.Range("J138").Value = 71990
$ws.Range("L138").Value = 71990
$ws.Range("N138").Value = -82270
